$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (ANDIKA): Lulus? flag changed from "v" to "x"
$ws.Range("F3").Value = "x"

# Row 7 (TEST123): Nilai changed 88.00 -> 0.00, Lulus? flag "v" -> "x", Keterangan "Lulus" cleared
$ws.Range("E7").Value = "0.00"
$ws.Range("F7").Value = "x"
$ws.Range("G7").Value = ""

# Row 8 (Ferlinyy): Lulus? flag changed from "v" to "x"
$ws.Range("F8").Value = "x"

# Row 9 (TITO): Lulus? flag changed from "v" to "x"
$ws.Range("F9").Value = "x"

# Row 10 (DITO): Nilai newly filled in as 0.00
$ws.Range("E10").Value = "0.00"

# Row 21 (HARRY POTTER): Nilai newly filled in as 90.00
$ws.Range("E21").Value = "90.00"
